$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$nombreLargo = "Conciencia histórica I. Perspectivas del México antiguo en los contextos globales"

# NC, Paterno, Materno, Nombres, Nombre_Largo, Grupo, Reprobadas
$data = @(
    @("23330051920005", "CASTRO",      "ARIAS",      "OMAR DAVID",       $nombreLargo, "4AEM",  3),
    @("23330051920018", "RAMOS",       "UTRERA",     "CARLOS DAVID",     $nombreLargo, "4AEM",  3),
    @("23330051920081", "CARRERA",     "MOLINA",     "MARIA DEL CARMEN", $nombreLargo, "4ALCM", 3),
    @("23330051920045", "SANTIAGO",    "GARCIA",     "URIEL",            $nombreLargo, "4BEM",  3),
    @("23330051920023", "VASQUEZ",     "ESPINDOLA",  "JOSUE YAHIR",      $nombreLargo, "4AEM",  2),
    @("23330051920025", "XOTLANIHUA",  "ZEPAHUA",    "JUAN ALBERTO",     $nombreLargo, "4AEM",  2),
    @("23330051920037", "HERNANDEZ",   "MARCELINO",  "LEONEL",           $nombreLargo, "4BEM",  2),
    @("22330051920021", "MEJIA",       "CRUZ",       "JOSE FRANCISCO",   $nombreLargo, "4BEM",  2),
    @("23330051920212", "VERA",        "VILLA",      "ALEX URIEL",       $nombreLargo, "4BEM",  2),
    @("22330061460232", "ALVAREZ",     "VOTE",       "CAMILO",           $nombreLargo, "4BLCM", 2),
    @("22330051920424", "COLMENARES",  "MARTINEZ",   "JULIO EDUARDO",    $nombreLargo, "4BLCM", 2),
    @("23330051920142", "REYES",       "ACEVEDO",    "KAREN AMERICA",    $nombreLargo, "4BLCM", 2),
    @("23330051920097", "MARIN",       "RODRIGUEZ",  "ABRIL",            $nombreLargo, "4ALCM", 1),
    @("23330051920133", "GUTIERREZ",   "PAZ",        "DANIA LIZETH",     $nombreLargo, "4BLCM", 1)
)

# Write column-by-column (A, then B, then C, then D, then E, then F, then G)
# so that new shared-string entries are interned in the same order the
# original authoring tool produced them (column-major scan).
for ($c = 0; $c -lt 7; $c++) {
    for ($i = 0; $i -lt $data.Count; $i++) {
        $row = $i + 2
        $vals = $data[$i]
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
